$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.07713289984429851
$ws.Range("C2").Value = 0.5125742293804644
$ws.Range("D2").Value = 0.5100403079960522
$ws.Range("E2").Value = 0.714171063538738
$ws.Range("F2").Value = 0.7195240788711913

# Row 3 (Q0)
$ws.Range("B3").Value = 0.1511691911167132
$ws.Range("C3").Value = 0.6823290450278029
$ws.Range("D3").Value = 0.8082812941544923
$ws.Range("E3").Value = 0.8990446563739158
$ws.Range("F3").Value = 0.889929442648896
$ws.Range("G3").Value = 121

# Row 4 (Q1)
$ws.Range("B4").Value = 0.09872457269350857
$ws.Range("C4").Value = 0.5659997234325101
$ws.Range("D4").Value = 0.4634181017742184
$ws.Range("E4").Value = 0.6807481926338244
$ws.Range("F4").Value = 0.6793331170920168
$ws.Range("G4").Value = 59
